$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new text in D2 explaining the change
$ws.Range("D2").Value = "In master branch it became RED"

# Make the cell font red for the Mouse row (A2:C2) and the new note cell (D2)
$ws.Range("A2:C2").Font.Color = 255

# Update the active selection to D2 to mirror the saved view state
$ws.Range("D2").Select()
